$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text type (source workbook stores these as
# inline strings, e.g. "0.860" / "233.00" / "12.20" which Excel would otherwise
# auto-convert to numbers and mangle -- so force Text format before writing).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.858.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.293.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.73"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.87"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.86"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.634.69"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.860"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.283.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.744.88"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.74"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.73%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.29"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.82"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0892"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.75%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.238"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.35"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.20"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.443"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.98%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.87%  "
